$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 386 (shifts existing rows 386:430 down to 387:431)
$ws.Rows("386:386").Insert()

# Populate the newly inserted row 386 with the new weekly data point
# (Macroferia Regional de Talca / Zapallo / Camote / 1a (guarda), week of 2023-05-31)
$ws.Range("A386").Value = 5
$ws.Range("B386").Value = "Macroferia Regional de Talca"
$ws.Range("C386").Value = "Maule"
$ws.Range("D386").Value = 45077
$ws.Range("E386").Value = 7
$ws.Range("F386").Value = 100112045
$ws.Range("G386").Value = "Zapallo"
$ws.Range("H386").Value = "Camote"
$ws.Range("I386").Value = "1a (guarda)"
$ws.Range("J386").Value = 800
$ws.Range("K386").Value = 250
$ws.Range("L386").Value = 280
$ws.Range("M386").Value = 265
$ws.Range("N386").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O386").Value = "Región del Maule"
$ws.Range("P386").Value = 265
$ws.Range("Q386").Value = 1
$ws.Range("R386").Value = "Hortaliza"
